# Apply odds updates to Sheet1 as described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 changes
$ws.Range("O2").Value = 1.25
$ws.Range("P2").Value = 3.75
$ws.Range("Q2").Value = 1.88
$ws.Range("R2").Value = 1.98

# Row 3 changes
$ws.Range("G3").Value = 2.4
$ws.Range("I3").Value = 3.25
$ws.Range("L3").Value = 4.33
$ws.Range("W3").Value = 5.5
$ws.Range("AC3").Value = 5.5
$ws.Range("AI3").Value = 15
$ws.Range("AZ3").Value = 81

# Row 5 changes
$ws.Range("M5").Value = 1.08
$ws.Range("N5").Value = 8
$ws.Range("Q5").Value = 2.3
$ws.Range("R5").Value = 1.6
$ws.Range("AG5").Value = 1000
